$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Ensure target cells keep a literal text format (they already store plain
# text such as "304.31" / "6.35%" as inline strings, not numbers/percentages)
$cells = @{
    'D2' = '304.31'
    'E2' = '6.35%'
    'D3' = '32.02'
    'E3' = '9.28%'
    'D4' = '5.286'
    'E4' = '4.30%'
    'D5' = '0.07465'
    'D6' = '7.851'
    'E6' = '7.02%'
    'D7' = '3.771'
    'E7' = '9.78%'
    'D8' = '1.468'
    'E8' = '6.32%'
    'D9' = '0.9154'
    'E9' = '1.58%'
    'D10' = '0.01737'
    'E10' = '2,583.57%'
    'D11' = '0.1698'
    'E11' = '6.90%'
    'D12' = '0.07708'
    'E12' = '11.94%'
    'D13' = '0.08013'
    'E13' = '5.42%'
    'D14' = '0.03034'
    'E14' = '3.76%'
    'D15' = '0.09846'
    'E15' = '9.55%'
    'D16' = '0.001498'
    'D17' = '0.04561'
    'E17' = '1.81%'
    'D18' = '0.006320'
    'E18' = '-4.32%'
    'D19' = '3.475'
    'E19' = '0.50%'
    'D20' = '2.228'
    'E20' = '-0.12%'
    'D21' = '0.3302'
    'E21' = '3.00%'
    'D22' = '0.1346'
    'E22' = '2.01%'
    'D23' = '4.484'
    'E23' = '12.04%'
    'D24' = '0.1642'
    'E24' = '3.99%'
    'E25' = '1.28%'
    'D26' = '0.004410'
    'E26' = '0.83%'
    'D27' = '0.0001400'
    'E27' = '19.90%'
    'D28' = '0.0001778'
    'E28' = '10.07%'
    'D40' = '0.04504'
    'E40' = '6.14%'
    'D41' = '0.007216'
    'E41' = '5.92%'
    'D42' = '0.1338'
    'E42' = '8.03%'
    'D43' = '0.002219'
    'E43' = '1.57%'
    'D44' = '0.01349'
    'E44' = '17.56%'
    'D45' = '0.00006225'
    'E45' = '8.53%'
    'D46' = '1.873'
    'E46' = '-2.78%'
    'D47' = '0.01302'
    'E47' = '-13.27%'
}

foreach ($addr in $cells.Keys) {
    $rng = $ws.Range($addr)
    $rng.NumberFormat = "@"
    $rng.Value = $cells[$addr]
}
